# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 15:05"

# --- Reorder Sierra Leona / Republica del Chad (rows 127 & 128) ---
# Row 127 becomes Sierra Leona (with refreshed stats), row 128 becomes
# Republica del Chad (keeping its previous stats, now shown on row 128).
$ws.Range("A127").Value = "Sierra Leona"
$ws.Range("B127").Value = 570
$ws.Range("C127").Value = 36
$ws.Range("D127").Value = 205
$ws.Range("E127").Value = 331
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 34

$ws.Range("A128").Value = "Republica del Chad"
$ws.Range("B128").Value = 545
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 139
$ws.Range("E128").Value = 350
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 56

# --- Update per-country statistic rows (Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 11 - Alemania
$ws.Range("B11").Value = 177910
$ws.Range("C11").Value = 83
$ws.Range("D11").Value = 156900
$ws.Range("E11").Value = 12810
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = 8200

# Row 18 - Arabia Saudita
$ws.Range("B18").Value = 62545
$ws.Range("C18").Value = 2691
$ws.Range("D18").Value = 33478
$ws.Range("E18").Value = 28728
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 339

# Row 29 - Portugal
$ws.Range("B29").Value = 29660
$ws.Range("C29").Value = 228
$ws.Range("D29").Value = 6452
$ws.Range("E29").Value = 21945
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = 1263

# Row 49 - Serbia
$ws.Range("B49").Value = 10833
$ws.Range("C49").Value = 100
$ws.Range("D49").Value = 5067
$ws.Range("E49").Value = 5531
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 235

# Row 51 - Argentina
$ws.Range("D51").Value = 2933
$ws.Range("E51").Value = 5482
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 394

# Row 75 - Uzbekistan
$ws.Range("B75").Value = 2927
$ws.Range("C75").Value = 72
$ws.Range("D75").Value = 2369
$ws.Range("E75").Value = 545

# Row 88 - Islandia
$ws.Range("B88").Value = 1803
$ws.Range("C88").Value = 1
$ws.Range("E88").Value = 4
